$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 62 ---
$ws.Range("A61:Z61").Copy()
$ws.Range("A62:Z62").PasteSpecial(-4122)
$ws.Range("AB61:AB61").Copy()
$ws.Range("AB62:AB62").PasteSpecial(-4122)
$ws.Range("A62").Value = 43545.438302905095
$ws.Range("B62").Value = "Entre 5 et 8 ans"
$ws.Range("C62").Value = "Conférences (Jancovici, Bihouix etc.), Articles de vulgarisation & blogs, Vidéos Youtube de vulgarisation"
$ws.Range("D62").Value = 6
$ws.Range("E62").Value = 8
$ws.Range("F62").Value = "Va commencer dans les 10 ans qui viennent"
$ws.Range("G62").Value = "Un peu plus lent (de l'ordre de 20 à 30 ans)"
$ws.Range("H62").Value = "Tristesse"
$ws.Range("I62").Value = "On a un fort potentiel d'action à l'échelle individuelle, Je comprend qu'il y ait des personnes climatosceptiques au sein de la population, Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres, Une transition écologique efficace peut se faire dans un cadre démocratique"
$ws.Range("J62").Value = "Problème d'éducation et/ou d'information"
$ws.Range("K62").Value = 8
$ws.Range("L62").Value = 5
$ws.Range("M62").Value = 7
$ws.Range("N62").Value = 5
$ws.Range("O62").Value = 6
$ws.Range("P62").Value = "Je partage des références directement à mon entourage (mail, vive voix etc.), Je donne une conférence sur les enjeux énergie-climat (hi hi hi)"
$ws.Range("Q62").Value = 3
$ws.Range("R62").Value = "En effet, on m'a déjà fait la remarque que dans mon argumentation, je souhaite presque cet effondrement. Tous mes arguments et les réponses vont vers ça. Comme pour le mouvement écolo des années 70 qui était déjà aussi alarmiste.
C'est vrai que pour ma part en tout cas, je sens que dans la crise écolo et sociale, il y a la confirmation à mon intuition/constat personnel que le système est mauvais et qu'il faut le changer. Intuition que j'ai depuis toujours ou presque. Et là ces constats alarmistes et scientifiques me donnent raison, donc j'y vais.
Et j'y vais d'autant mieux qu'on se retrouve, qu'on construit/réfléchis ensemble, bref qu'on retrouve un lien social fort autour d'une thématique qui est très forte chez nous, et donc qui recréé un sentiment d'appartenance à un groupe, une communauté, etc. Ce qui est un sentiment génial et rassurant, d'autant plus en ces temps où la famille est moins forte/nombreuse, les amitiés moins durables ou solides...
Et dans ce confortable sentiment d'appartenance, où on se reconnaît entre nous, certains peuvent devenir dogmatique et `"dériver`"vers une opposition ou un besoin de se définir contre les autres."
$ws.Range("S62").Value = "Une femme"
$ws.Range("T62").Value = 32
$ws.Range("U62").Value = "En ville dans une grande agglomération"
$ws.Range("V62").Value = "Études supérieures longues (ingénieur, école de commerce, Master MBA graduate à l'étranger)"
$ws.Range("W62").Value = "Langues / Littérature / Communication"
$ws.Range("X62").Value = "Très frugal (flexitarien ou végétalien, AMAP, déplacement doux)"
$ws.Range("Y62").Value = "Local, Bio majoritairement (+ de 50% de ce que tu manges chez toi), Flexitarien"
$ws.Range("Z62").Value = "Très à gauche (France insoumise ou plus à gauche en France)"
$ws.Range("AB62").Value = "Académicien"
$ws.Rows(62).EntireRow.AutoFit()

# --- Row 63 ---
$ws.Range("A61:Q61").Copy()
$ws.Range("A63:Q63").PasteSpecial(-4122)
$ws.Range("S61:Z61").Copy()
$ws.Range("S63:Z63").PasteSpecial(-4122)
$ws.Range("AB61:AB61").Copy()
$ws.Range("AB63:AB63").PasteSpecial(-4122)
$ws.Range("A63").Value = 43545.490691018524
$ws.Range("B63").Value = "Entre 2 et 5 ans"
$ws.Range("C63").Value = "Conférences (Jancovici, Bihouix etc.), Livres, Cours"
$ws.Range("D63").Value = 5
$ws.Range("E63").Value = 6
$ws.Range("F63").Value = "A déjà commencé"
$ws.Range("G63").Value = "Un peu plus lent (de l'ordre de 20 à 30 ans)"
$ws.Range("H63").Value = "Peur"
$ws.Range("I63").Value = "Je suis prêt à baisser mon niveau de vie si cette baisse s'opère pour les autres également, Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres"
$ws.Range("J63").Value = "Une réaction de personnes ne voulant pas toucher à leur mode de vie"
$ws.Range("K63").Value = 6
$ws.Range("L63").Value = 7
$ws.Range("M63").Value = 7
$ws.Range("N63").Value = 5
$ws.Range("O63").Value = 5
$ws.Range("P63").Value = "Je partage des liens sur les réseaux sociaux, Je partage des références directement à mon entourage (mail, vive voix etc.), Je donne une conférence sur les enjeux énergie-climat (hi hi hi)"
$ws.Range("Q63").Value = 2
$ws.Range("S63").Value = "Une femme"
$ws.Range("T63").Value = 22
$ws.Range("U63").Value = "En ville dans une grande agglomération"
$ws.Range("V63").Value = "Études supérieures longues (ingénieur, école de commerce, Master MBA graduate à l'étranger)"
$ws.Range("W63").Value = "Sciences sociales"
$ws.Range("X63").Value = "Plutôt écolo (vélo, transport en commun, limitation de la consommation et notamment de la viande)"
$ws.Range("Y63").Value = "Flexitarien"
$ws.Range("Z63").Value = "Très à gauche (France insoumise ou plus à gauche en France)"
$ws.Range("AB63").Value = "Académicien"

# --- Row 64 ---
$ws.Range("A61:Q61").Copy()
$ws.Range("A64:Q64").PasteSpecial(-4122)
$ws.Range("S61:Z61").Copy()
$ws.Range("S64:Z64").PasteSpecial(-4122)
$ws.Range("AB61:AB61").Copy()
$ws.Range("AB64:AB64").PasteSpecial(-4122)
$ws.Range("A64").Value = 43545.62710184028
$ws.Range("B64").Value = "Entre 2 et 5 ans"
$ws.Range("C64").Value = "Conférences (Jancovici, Bihouix etc.), Vidéos Youtube de vulgarisation"
$ws.Range("D64").Value = 4
$ws.Range("E64").Value = 8
$ws.Range("F64").Value = "A déjà commencé"
$ws.Range("G64").Value = "Plutôt rapide (5 à 10 ans)"
$ws.Range("H64").Value = "Colère"
$ws.Range("I64").Value = "Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres, Une transition écologique efficace peut se faire dans un cadre démocratique"
$ws.Range("J64").Value = "Problème d'éducation et/ou d'information"
$ws.Range("K64").Value = 6
$ws.Range("L64").Value = 7
$ws.Range("M64").Value = 6
$ws.Range("N64").Value = 5
$ws.Range("O64").Value = 5
$ws.Range("P64").Value = "Je partage des liens sur les réseaux sociaux, Je partage des références directement à mon entourage (mail, vive voix etc.), Je donne une conférence sur les enjeux énergie-climat (hi hi hi)"
$ws.Range("Q64").Value = 3
$ws.Range("S64").Value = "Une femme"
$ws.Range("T64").Value = 29
$ws.Range("U64").Value = "En ville dans une grande agglomération"
$ws.Range("V64").Value = "Études supérieures longues (ingénieur, école de commerce, Master MBA graduate à l'étranger)"
$ws.Range("W64").Value = "Energie"
$ws.Range("X64").Value = "Plutôt écolo (vélo, transport en commun, limitation de la consommation et notamment de la viande)"
$ws.Range("Y64").Value = "Bio majoritairement (+ de 50% de ce que tu manges chez toi), Flexitarien"
$ws.Range("Z64").Value = "Très à gauche (France insoumise ou plus à gauche en France)"
$ws.Range("AB64").Value = "Académicien"

# --- Row 65 ---
$ws.Range("A61:Z61").Copy()
$ws.Range("A65:Z65").PasteSpecial(-4122)
$ws.Range("AB61:AB61").Copy()
$ws.Range("AB65:AB65").PasteSpecial(-4122)
$ws.Range("A65").Value = 43547.43574328704
$ws.Range("B65").Value = "Depuis moins de 2 ans"
$ws.Range("C65").Value = "Conférences (Jancovici, Bihouix etc.), Articles de vulgarisation & blogs, Livres"
$ws.Range("D65").Value = 4
$ws.Range("E65").Value = 7
$ws.Range("F65").Value = "Va commencer dans les 15 à 20 ans qui viennent"
$ws.Range("G65").Value = "Un peu plus lent (de l'ordre de 20 à 30 ans)"
$ws.Range("H65").Value = "Colère"
$ws.Range("I65").Value = "On a un fort potentiel d'action à l'échelle individuelle, Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres, Une transition écologique efficace peut se faire dans un cadre démocratique"
$ws.Range("J65").Value = "Une réaction de personnes ne voulant pas toucher à leur mode de vie, Problème d'éducation et/ou d'information, Des personnes constituées d'un système de valeurs remis en question quand il s'agit de parler de transition écologique, énergétique, de sobriété, de low tech"
$ws.Range("K65").Value = 7
$ws.Range("L65").Value = 6
$ws.Range("M65").Value = 7
$ws.Range("N65").Value = 7
$ws.Range("O65").Value = 5
$ws.Range("P65").Value = "Je partage des liens sur les réseaux sociaux, Je partage des références directement à mon entourage (mail, vive voix etc.), Je fais appel à l'émotion (par exemple à des parents s'inquiétant pour leurs enfants)., Je donne une conférence sur les enjeux énergie-climat (hi hi hi)"
$ws.Range("Q65").Value = 2
$ws.Range("R65").Value = "Il est vrai que l'on fait confiance aux scientifiques et que l'on partage facilement des articles sur le sujet du dérèglement climatique parfois aveuglement. cependant chacun a son échelle et même un européen voit un peu les impacts nocifs : pollution, vague de chaleur, inondation, feu de foret, perte de biodiversité, artificialisation massive, pollution plastique dans les cous d'eau.. donc même si on n'est pas scientifique et que l'on n'a pas tous accès à l'information exacte, on peut faire état de la situation sans l'avoir entendu par un gourou"
$ws.Range("S65").Value = "Une femme"
$ws.Range("T65").Value = 23
$ws.Range("U65").Value = "En banlieue d'un grand centre urbain"
$ws.Range("V65").Value = "Études supérieures longues (ingénieur, école de commerce, Master MBA graduate à l'étranger)"
$ws.Range("W65").Value = "Economie / Gestion"
$ws.Range("X65").Value = "Plutôt écolo (vélo, transport en commun, limitation de la consommation et notamment de la viande)"
$ws.Range("Y65").Value = "Local, Bio majoritairement (+ de 50% de ce que tu manges chez toi), Flexitarien"
$ws.Range("Z65").Value = "Parti à préoccupation environnementale (Europe Écologie les Verts en France)"
$ws.Range("AB65").Value = "Académicien"

# --- Row 66 ---
$ws.Range("A61:Z61").Copy()
$ws.Range("A66:Z66").PasteSpecial(-4122)
$ws.Range("AB61:AB61").Copy()
$ws.Range("AB66:AB66").PasteSpecial(-4122)
$ws.Range("A66").Value = 43547.54042091435
$ws.Range("B66").Value = "Entre 5 et 8 ans"
$ws.Range("C66").Value = "Articles de vulgarisation & blogs, Livres, Vidéos Youtube de vulgarisation, Articles de presse"
$ws.Range("D66").Value = 5
$ws.Range("E66").Value = 10
$ws.Range("F66").Value = "A déjà commencé"
$ws.Range("G66").Value = "Plutôt rapide (5 à 10 ans)"
$ws.Range("H66").Value = "Soulagement"
$ws.Range("I66").Value = "On peut augmenter le pouvoir d'achat en France tout en se limitant à une augmentation de la T° moyenne de 2°C, Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres"
$ws.Range("J66").Value = "Problème d'éducation et/ou d'information"
$ws.Range("K66").Value = 5
$ws.Range("L66").Value = 7
$ws.Range("M66").Value = 7
$ws.Range("N66").Value = 7
$ws.Range("O66").Value = 7
$ws.Range("P66").Value = "Je partage des références directement à mon entourage (mail, vive voix etc.), Je fais profil bas. Trop en parler, c'est devenir prêcheur, et donc desservir la cause., Je donne une conférence sur les enjeux énergie-climat (hi hi hi)"
$ws.Range("Q66").Value = 3
$ws.Range("R66").Value = "Je ressens profondément le fait d'être une bulle, une bulle confortable et consensuelle, vis à vis de mon militantisme écologique. Je pense que j'ai une forme de fracture dans mon esprit que les contestations sociales représentées par le mouvement des gilets jaunes contribuent à rendre plus intelligibles, que j'arrive de mieux à mieux à appréhender et à verbaliser. Ça me donne le sentiment que je me suis trompée de combat avec le militantisme écologique et que la crise sociale est bien plus urgente pour l'humanité."
$ws.Range("S66").Value = "Une femme"
$ws.Range("T66").Value = 26
$ws.Range("U66").Value = "En banlieue d'un grand centre urbain"
$ws.Range("V66").Value = "Études supérieures courtes (DUT BTS ou licence pro en France, Bachelor à l'étranger)"
$ws.Range("W66").Value = "Economie / Gestion"
$ws.Range("X66").Value = "Plutôt écolo (vélo, transport en commun, limitation de la consommation et notamment de la viande)"
$ws.Range("Y66").Value = "Local, Bio majoritairement (+ de 50% de ce que tu manges chez toi), Flexitarien"
$ws.Range("Z66").Value = "Aucun de ces partis, je ne crois pas à la politique"
$ws.Range("AB66").Value = "Académicien"

$excel.CutCopyMode = 0
Write-Host "Rows 62-66 added"
